# Super Herois workbook update:
#  - append new heroes to the bottom of the list
#  - rename the header cell (A1) from "nome" to "target"
#  - highlight duplicate values in column A (conditional formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new hero names below the existing list (rows 16-26).
$ws.Cells.Item(16, 1).Value = "Flash "
$ws.Cells.Item(17, 1).Value = "Spider Man"
$ws.Cells.Item(18, 1).Value = "Hulk"
$ws.Cells.Item(19, 1).Value = "Wolverine"
$ws.Cells.Item(20, 1).Value = "Deadpool"
$ws.Cells.Item(21, 1).Value = "Pantera Negra"
$ws.Cells.Item(22, 1).Value = "Adão Negro"
$ws.Cells.Item(23, 1).Value = "Demolidor"
$ws.Cells.Item(24, 1).Value = "Doutor Strange"
$ws.Cells.Item(25, 1).Value = "Viúva negra"
$ws.Cells.Item(26, 1).Value = "Green Arrow"

# Rename the header from "nome" to "target".
$ws.Range("A1").Value = "target"

# Highlight duplicate values across the whole column (Excel's built-in
# "Light Red Fill with Dark Red Text" duplicate-values preset).
$rng = $ws.Range("A1:A1048576")
$fc = $rng.FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1  # xlDuplicate (0 = xlUnique)
$fc.Font.Color = 393372      # RGB(156,0,6)   -> dark red text
$fc.Interior.Color = 13551615  # RGB(255,199,206) -> light red fill

# Adjust the view: jump to the top, select A2, set zoom to 130%.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 130
